# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated handoff timestamps across the
# Overview / zh-cn / de-de sheets. Also widens the "Status" columns that
# now need to fit the longer text.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Sheets.Item("Overview")
$ws_zhcn     = $wb.Sheets.Item("zh-cn")
$ws_dede     = $wb.Sheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
# E2/F2 hold the per-language status, G2 the "Latest HO Xliff Generate
# Date" timestamp for the row.
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_overview.Range("G2").Value = "2016-08-27 04:57:17"

# --- zh-cn sheet ------------------------------------------------------
# C2 is Status, H2 is Latest Handoff Datetime.
$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("H2").Value = "2016-08-27 04:57:13"

# --- de-de sheet ------------------------------------------------------
# C2 is Status, H2 is Latest Handoff Datetime.
$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("H2").Value = "2016-08-27 04:57:17"

# --- Column widths ------------------------------------------------------
# The Status columns grew to fit "Ready for handoff". ColumnWidth is
# expressed in characters and Excel snaps it to whole-pixel increments,
# so we target the closest attainable width to the authored value
# (~17.22 characters).
$newStatusColWidth = 16.3333333333333

$ws_overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$ws_overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$ws_zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$ws_dede.Columns.Item(3).ColumnWidth = $newStatusColWidth
